$ws = $excel.ActiveWorkbook.ActiveSheet
$ws.Range("C59").Value = "/*/*/oos:protocolNumber"
Write-Output $ws.Range("C59").Value()
